# Applies the "Atualização de bases das ligas" edit to the "Turkey 1 Lig" sheet.
# The edit swaps the B:AC content between specific pairs of adjacent rows
# (column A, the sequential index, stays put), and applies a handful of
# single-cell corrections on a few other rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Turkey 1 Lig")

function Swap-RowRange {
    param(
        $Row1,
        $Row2
    )

    $rng1 = $ws.Range($ws.Cells.Item($Row1, 2), $ws.Cells.Item($Row1, 29))
    $rng2 = $ws.Range($ws.Cells.Item($Row2, 2), $ws.Cells.Item($Row2, 29))

    $vals1 = $rng1.Value()
    $vals2 = $rng2.Value()

    $rng1.Value = $vals2
    $rng2.Value = $vals1
}

# Row pairs whose B:AC contents were swapped
$pairs = @(
    @(32, 33),
    @(37, 38),
    @(43, 44),
    @(48, 49),
    @(85, 86),
    @(90, 91),
    @(157, 158),
    @(196, 197),
    @(232, 233)
)

foreach ($pair in $pairs) {
    Swap-RowRange $pair[0] $pair[1]
}

# Individual cell corrections on rows without a swap partner
$ws.Range("R246").Value = 1.775
$ws.Range("S246").Value = 2.025

$ws.Range("N248").Value = 2.3
$ws.Range("O248").Value = 3.25
$ws.Range("P248").Value = 2.8
$ws.Range("Q248").Value = -0.25
$ws.Range("R248").Value = 2.025
$ws.Range("S248").Value = 1.775
$ws.Range("U248").Value = 1.95
$ws.Range("V248").Value = 1.85

$ws.Range("R251").Value = 1.9
$ws.Range("S251").Value = 1.9

$ws.Range("N252").Value = 8
$ws.Range("O252").Value = 5
$ws.Range("R252").Value = 1.825
$ws.Range("S252").Value = 1.975
$ws.Range("T252").Value = 2.5
